$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the contents of rows 18 and 19 (columns D, E, F, G, I) ---
$d18 = $ws.Range("D18").Value
$e18 = $ws.Range("E18").Value
$f18 = $ws.Range("F18").Value
$g18 = $ws.Range("G18").Value
$i18 = $ws.Range("I18").Value

$d19 = $ws.Range("D19").Value
$e19 = $ws.Range("E19").Value
$f19 = $ws.Range("F19").Value
$g19 = $ws.Range("G19").Value
$i19 = $ws.Range("I19").Value

$ws.Range("D18").Value = $d19
$ws.Range("E18").Value = $e19
$ws.Range("F18").Value = $f19
$ws.Range("G18").Value = $g19
$ws.Range("I18").Value = $i19

$ws.Range("D19").Value = $d18
$ws.Range("E19").Value = $e18
$ws.Range("F19").Value = $f18
$ws.Range("G19").Value = $g18
$ws.Range("I19").Value = $i18

# --- Insert two new rows (new sightings at Mont-Tremblant) before row 28 ---
$ws.Rows("28:29").Insert()

$ws.Range("A28").Value = "5/3/2022"
$ws.Range("B28").Value = "RASY"
$ws.Range("C28").Value = 153
$ws.Range("D28").Value = "Mont-Tremblant"
$ws.Range("E28").Value = "Laurentides"
$ws.Range("F28").Value = "C"
$ws.Range("G28").Value = "Cote 1"
$ws.Range("H28").Value = ""
$ws.Range("I28").Value = "Jacques Tremblay"

$ws.Range("A29").Value = "5/3/2022"
$ws.Range("B29").Value = "PSTR"
$ws.Range("C29").Value = 153
$ws.Range("D29").Value = "Mont-Tremblant"
$ws.Range("E29").Value = "Laurentides"
$ws.Range("F29").Value = "C"
$ws.Range("G29").Value = "Cote 1"
$ws.Range("H29").Value = ""
$ws.Range("I29").Value = "Jacques Tremblay"

$ws.Range("C35").Select()
